# "9th Stab- Cosmetic Changes"
# Adds two new snapshot columns (Jun_15 / Jun_17) to the ratings grid.
#
# Before:  B1="Jun_13"  C1="Jun_10"
# After:   B1="Jun_17"  C1="Jun_15"  D1="Jun_13"  E1="Jun_10"
#
# For every data row (2-27) the "UN" ticker value that lived in column B
# is echoed into two brand-new columns (C and D) while the original
# rating-action text that used to live in column C slides over to the
# new column E - exactly as if two blank columns had been inserted
# immediately to the right of column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert two new (blank) columns at C/D. This pushes the old
# column C (rating-action text) out to column E on every row, while
# column A/B are left untouched.
$ws.Range("C1:D1").EntireColumn.Insert()

# Step 2: fix up the header row. The old B1 value ("Jun_13") needs to
# move to D1, and the new B1/C1 cells get the two newest week labels.
$oldB1 = $ws.Range("B1").Value2
$ws.Range("D1").Value = $oldB1
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# Step 3: populate the two new data columns (C, D) with the same "UN"
# value that's already sitting in column B for every analyst row.
for ($r = 2; $r -le 27; $r++) {
    $unVal = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 3).Value = $unVal
    $ws.Cells.Item($r, 4).Value = $unVal
}

# Step 4: keep the custom 8-character width that column C originally
# had, now mirrored across the three data columns C, D and E.
$ws.Columns.Item(3).ColumnWidth = 7.083333333333333
$ws.Columns.Item(4).ColumnWidth = 7.083333333333333
$ws.Columns.Item(5).ColumnWidth = 7.083333333333333
